$d = $word.ActiveDocument

# Locate the paragraph that holds the screenshot (drawing) so the new
# content is anchored after it regardless of exact paragraph index.
$drawingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.InlineShapes.Count -gt 0) {
        $drawingIndex = $i
    }
}

$afterDrawing = $d.Paragraphs.Item($drawingIndex + 1)

# Insert a brand new blank paragraph right before the paragraph that
# immediately follows the drawing (this becomes the blank line called
# for by the diff, while the pre-existing blank paragraphs are left
# alone and simply shift down).
$insertPos = $d.Range($afterDrawing.Range.Start, $afterDrawing.Range.Start)
$insertPos.InsertParagraphBefore()

# The paragraph that used to follow the drawing directly is now one
# slot further down; reuse it (instead of minting yet another new
# paragraph) so its formatting/identity stays untouched and becomes
# the "Code available" line.
$codeLinePara = $d.Paragraphs.Item($drawingIndex + 2)
$codeRange = $d.Range($codeLinePara.Range.Start, $codeLinePara.Range.Start)
$codeRange.InsertAfter("Code available: ")

$urlInsertPos = $d.Range($codeLinePara.Range.End - 1, $codeLinePara.Range.End - 1)
$urlInsertPos.InsertAfter("https://github.com/michaelmccoll/PDA_repository/tree/main/IT5")

# Nudge formatting off/back-on so the URL text is kept in its own run
# (matching the two separate <w:r> elements in the target revision)
# rather than being silently coalesced into the preceding run.
$urlInsertPos.Bold = 1
$urlInsertPos.Bold = 0
